$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Source sheet whose layout/style (A1:B11, header row + 10 data rows) is
# reused as the template for the two new plot-data sheets.
# ---------------------------------------------------------------------------
$extreme = $wb.Worksheets.Item("extreme_vols")

# ===========================================================================
# 1) annual_vols  (sheetId 6) - Tenors 1..10, Vols 0.1..1.0
# ===========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extreme.Copy($null, $lastSheet)
$annual = $wb.Worksheets.Item($wb.Worksheets.Count)
$annual.Name = "annual_vols"

# The copied sheet inherits B10/B11's "table-bottom" styling (s=22 / s=5);
# normalize them back to the plain numeric style used by the rest of col B.
$annual.Range("B2").Copy()
$annual.Range("B10:B11").PasteSpecial(-4122)

# Column A: integer tenors 1..10. Column B: vols 0.1..1.0
for ($i = 0; $i -lt 10; $i++) {
    $annual.Cells.Item($i + 2, 1).Value = $i + 1
    $annual.Cells.Item($i + 2, 2).Value = [math]::Round(($i + 1) * 0.1, 1)
}

$annual.Range("A1:B11").Select()

# ===========================================================================
# 2) stepwise_vols  (sheetId 7) - Tenors 0.1..1.0, Vols 0.1..1.0
# ===========================================================================
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$extreme.Copy($null, $lastSheet2)
$stepwise = $wb.Worksheets.Item($wb.Worksheets.Count)
$stepwise.Name = "stepwise_vols"

# This sheet has no thick bottom border under the last two rows, so drop the
# special row height/border flag the template carried on rows 10-11.
$stepwise.Rows("10:11").AutoFit()

# Normalize A9:A11 / B10:B11 styling back to the plain numeric style (s=6/s=4)
$stepwise.Range("A2").Copy()
$stepwise.Range("A9:A11").PasteSpecial(-4122)
$stepwise.Range("B2").Copy()
$stepwise.Range("B10:B11").PasteSpecial(-4122)

# Both columns: 0.1 .. 1.0
for ($i = 0; $i -lt 10; $i++) {
    $v = [math]::Round(($i + 1) * 0.1, 1)
    $stepwise.Cells.Item($i + 2, 1).Value = $v
    $stepwise.Cells.Item($i + 2, 2).Value = $v
}

# ===========================================================================
# 3) extreme_vols selection tidy-up - now that it's no longer the active tab,
#    its lingering single-cell selection is replaced by a full-range one.
# ===========================================================================
$extreme.Range("A1:B11").Select()

# stepwise_vols ends up as the active tab/selection
$stepwise.Range("A2:A11").Select()
$stepwise.Activate()
